$wb = $excel.ActiveWorkbook

# Update status text from "Ready for handoff" to "In Translation" across all sheets,
# then resize the affected columns to fit the new (shorter) text.
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Columns("E:F").AutoFit()

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "In Translation"
$ws2.Columns("C:C").AutoFit()

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "In Translation"
$ws3.Columns("C:C").AutoFit()
